$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDTSlot")
$ws.Range("D1").Value = "dt"
$ws.Range("D1").Font.Name = "Calibri"
